# Replace the per-row price/volume/ticker-ref data (columns D-I) for rows 2-42
# with corrected values, and normalize the "fixed_ticker" column (I) to "ASML"
# for every row (it previously held a mix of unrelated tickers left over from
# other companies' extra data files).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: D (open_price), E (close_price), F (high_price), G (low_price), H (shares_outstanding)
$rowData = @{
  2 = @(90.92885540984577, 96.25222015380859, 99.9929605364234, 84.48146218578815, 388147674)
  3 = @(97.60996421740757, 89.84574127197266, 97.99953462054374, 87.68045670237217, 388147674)
  4 = @(78.13148860854035, 84.06564331054688, 85.53332917877151, 76.5097882420616, 388147674)
  5 = @(79.08275603037816, 83.20494842529297, 83.28648997314409, 69.91426365048702, 388147674)
  6 = @(90.05413391194676, 87.51739501953125, 92.95325815008673, 86.18560744661599, 388147674)
  7 = @(90.18045633534425, 100.4938659667969, 101.6214615030677, 85.11083665913154, 388147674)
  8 = @(100.2555232810406, 96.82688903808594, 100.7047276698905, 91.71143548688508, 388147674)
  9 = @(101.7406372699302, 111.2931442260742, 115.4001689737709, 100.4938601515366, 388147674)
  10 = @(121.3315197410946, 120.8731460571289, 123.7150545046432, 115.5468390933732, 388147674)
  11 = @(121.3788315615825, 139.1829071044922, 144.1824869601744, 119.7400745787592, 388147674)
  12 = @(157.7554746759066, 167.3472747802734, 167.5509627732638, 155.2649375930761, 388147674)
  13 = @(161.2274139085169, 187.910400390625, 193.113669496508, 161.0792746139941, 388147674)
  14 = @(182.8274679978861, 174.4763031005859, 197.529959852706, 173.7171139975621, 388147674)
  15 = @(180.8507457923005, 199.8866882324219, 207.0415142318016, 177.2920177822027, 388147674)
  16 = @(178.2821023923782, 160.9928436279297, 179.8232933008194, 147.5331809565631, 388147674)
  17 = @(144.105206198179, 163.4867401123047, 165.8779053895078, 137.6602680915779, 388147674)
  18 = @(178.6277101746051, 195.0483093261719, 195.2257668372502, 178.1793559553855, 388147674)
  19 = @(204.0882224516673, 210.5519104003907, 221.5987771190891, 193.1169561932205, 388147674)
  20 = @(235.4616662731024, 247.5574645996093, 254.5692587092618, 225.1802355326746, 388147674)
  21 = @(285.9316235671533, 266.3709411621094, 290.3259039851951, 265.9058974858189, 388147674)
  22 = @(240.2425767403694, 273.745361328125, 289.24397985695, 229.9164978181778, 388147674)
  23 = @(351.061438603179, 337.4112243652344, 384.2950860663689, 334.5686111342923, 388147674)
  24 = @(358.2156084474718, 344.5558471679688, 390.2473378984965, 340.9024480392986, 388147674)
  25 = @(478.7874784014135, 511.5084533691406, 549.4674986060739, 468.9721385956826, 388147674)
  26 = @(593.5170319524951, 620.616943359375, 646.9987173714082, 586.7947586174027, 388147674)
  27 = @(660.0014242670134, 736.3504638671875, 738.2328009138233, 638.4796641783259, 388147674)
  28 = @(715.4721983993965, 780.6618041992188, 789.3434860149313, 680.9278842848639, 388147674)
  29 = @(766.178561700371, 652.0422973632812, 769.6640754187806, 597.6123364827146, 388147674)
  30 = @(650.2032782788768, 542.8262329101562, 656.221093767734, 523.7906601501297, 388147674)
  31 = @(441.3986589690715, 556.9010620117188, 559.3053231198335, 400.0702732427433, 388147674)
  32 = @(411.6506633008337, 459.093505859375, 482.2026940787879, 352.9058869480995, 388147674)
  33 = @(547.0895063465185, 644.1323852539062, 667.6230548984161, 531.9716273143085, 388147674)
  34 = @(659.1669451368648, 622.1058959960938, 665.8093981029938, 591.726411301878, 388147674)
  35 = @(719.0835849546435, 701.860595703125, 756.3020466960675, 659.6653453680987, 388147674)
  36 = @(572.7863844551795, 587.9563598632812, 612.8763434057247, 557.2826052244899, 388147674)
  37 = @(718.6045542630967, 856.241943359375, 870.2005835185993, 685.2140690579396, 388147674)
  38 = @(966.1536293233318, 860.3535766601562, 1008.457821067475, 843.7080524005079, 388147674)
  39 = @(1023.292591916819, 925.6126098632812, 1096.950205034136, 840.1166986429427, 388147674)
  40 = @(831.2783607476397, 665.8206176757812, 864.9084966854472, 660.5439314288184, 388147674)
  41 = @(699.9033757713656, 733.6115112304688, 778.7707369281133, 661.4619267840154, 388147674)
  42 = @(660.841955796277, 664.332275390625, 688.4064252309196, 575.264730709894, 388147674)
}

foreach ($r in $rowData.Keys) {
  $vals = $rowData[$r]
  $ws.Range("D$r").Value = $vals[0]
  $ws.Range("E$r").Value = $vals[1]
  $ws.Range("F$r").Value = $vals[2]
  $ws.Range("G$r").Value = $vals[3]
  $ws.Range("H$r").Value = $vals[4]
  $ws.Range("I$r").Value = "ASML"
}

